$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aggregate c_DK_Central and c_DK_Decentral loads into B2 (c_DK_Central row)
$ws.Range("B2").Value2 = 38775999.99999999

# Remove row 3 (the now-aggregated c_DK_Decentral row)
$ws.Rows.Item(3).Delete()
